$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.147.01"
$ws.Range("E2").Value = "'  +0.31%  "
$ws.Range("D3").Value = "'2.116.63"
$ws.Range("E3").Value = "'  +0.66%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "'  +0.24%  "
$ws.Range("D5").Value = "'346.31"
$ws.Range("E5").Value = "'  +0.46%  "
$ws.Range("E6").Value = "'  -0.21%  "
$ws.Range("D7").Value = "'0.5211"
$ws.Range("E7").Value = "'  +0.55%  "
$ws.Range("D8").Value = "'0.4468"
$ws.Range("E8").Value = "'  -0.18%  "
$ws.Range("D9").Value = "'54.18"
$ws.Range("E9").Value = "'  +3.93%  "
$ws.Range("D10").Value = "'0.09354"
$ws.Range("E10").Value = "'  -1.46%  "
$ws.Range("D11").Value = "'1.184"
$ws.Range("D12").Value = "'25.36"
$ws.Range("E12").Value = "'  +0.75%  "
$ws.Range("D13").Value = "'8.659"
$ws.Range("E13").Value = "'  +7.34%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.976"
$ws.Range("E14").Value = "'  +3.48%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'2.103.63"
$ws.Range("E15").Value = "'  -0.37%  "
$ws.Range("D16").Value = "'102.74"
$ws.Range("E16").Value = "'  +3.40%  "
$ws.Range("D17").Value = "'0.00001173"
$ws.Range("E17").Value = "'  +0.28%  "
$ws.Range("D18").Value = "'1.008"
$ws.Range("E18").Value = "'  -0.15%  "
$ws.Range("D19").Value = "'21.57"
$ws.Range("E19").Value = "'  +4.77%  "
$ws.Range("D20").Value = "'0.06699"
$ws.Range("E20").Value = "'  -0.17%  "
$ws.Range("E21").Value = "'  +2.00%  "
$ws.Range("E22").Value = "'  -0.15%  "
$ws.Range("D23").Value = "'30.164.76"
$ws.Range("E23").Value = "'  +0.07%  "
$ws.Range("D24").Value = "'12.76"
$ws.Range("E24").Value = "'  +0.63%  "
$ws.Range("E25").Value = "'  +0.54%  "
$ws.Range("D26").Value = "'2.381.03"
$ws.Range("E26").Value = "'  +0.97%  "
$ws.Range("D27").Value = "'22.19"
$ws.Range("E27").Value = "'  +0.80%  "
$ws.Range("D28").Value = "'2.559"
$ws.Range("E28").Value = "'  +1.07%  "
$ws.Range("D29").Value = "'163.09"
$ws.Range("E29").Value = "'  -0.78%  "
$ws.Range("D30").Value = "'134.29"
$ws.Range("E30").Value = "'  +0.57%  "
$ws.Range("D31").Value = "'1.159"
$ws.Range("E31").Value = "'  +0.01%  "
$ws.Range("D32").Value = "'1.793"
$ws.Range("E32").Value = "'  +10.52%  "
$ws.Range("E33").Value = "'  +0.26%  "
$ws.Range("D34").Value = "'6.858"
$ws.Range("E34").Value = "'  +10.94%  "
$ws.Range("D35").Value = "'6.303"
$ws.Range("E35").Value = "'  +0.86%  "
$ws.Range("D36").Value = "'3.967"
$ws.Range("E36").Value = "'  +0.37%  "
$ws.Range("E37").Value = "'  +6.33%  "
$ws.Range("D38").Value = "'0.02651"
$ws.Range("E38").Value = "'  +2.97%  "
$ws.Range("D39").Value = "'0.06874"
$ws.Range("E39").Value = "'  +1.36%  "
$ws.Range("D40").Value = "'0.7166"
$ws.Range("E40").Value = "'  +3.18%  "
$ws.Range("D41").Value = "'12.80"
$ws.Range("E41").Value = "'  +2.87%  "
$ws.Range("E42").Value = "'  -1.31%  "
$ws.Range("D43").Value = "'1.335"
$ws.Range("E43").Value = "'  +1.78%  "
$ws.Range("D44").Value = "'0.6978"
$ws.Range("E44").Value = "'  +4.03%  "
$ws.Range("D45").Value = "'14.79"
$ws.Range("E45").Value = "'  +4.06%  "
$ws.Range("D46").Value = "'2.400"
$ws.Range("E46").Value = "'  +5.38%  "
$ws.Range("E47").Value = "'  -0.09%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.637"
$ws.Range("E48").Value = "'  -0.12%  "
$ws.Range("B49").Value = "WEMIXTOKEN"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'1.264"
$ws.Range("E49").Value = "'  +7.84%  "
$ws.Range("D50").Value = "'0.00000000350"
$ws.Range("E50").Value = "'  +3.30%  "
$ws.Range("D51").Value = "'1.218"
$ws.Range("E51").Value = "'  +9.28%  "
